$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.517.82'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '2.369.41'
$ws.Range('E3').Value = '  -4.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '85.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.533'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.43%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0828'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '30.22'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.96%  '
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '2.738.91'
$ws.Range('E13').Value = '  -4.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.88'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.48%  '
$ws.Range('D16').Value = '2.376.63'
$ws.Range('E16').Value = '  -4.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.756'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.10%  '
$ws.Range('D18').Value = '40.496.63'
$ws.Range('E18').Value = '  -2.76%  '
$ws.Range('D19').Value = '0.0₃0908'
$ws.Range('E19').Value = '  -3.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.57'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.32%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  -7.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.21'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.38'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.19'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0727'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.10%  '
$ws.Range('E35').Value = '  -5.39%  '
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.84'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0995'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.30%  '
$ws.Range('E40').Value = '  -8.07%  '
$ws.Range('E41').Value = '  -5.20%  '
$ws.Range('E42').Value = '  -3.98%  '
$ws.Range('D43').Value = '1.966.62'
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0267'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.92%  '
$ws.Range('D48').Value = '2.603.81'
$ws.Range('E48').Value = '  -4.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '92.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.09%  '
